# Apply new metrics values (all rows share the same B:Q values) and
# reorder the model names in column A according to the target permutation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared B:Q values (row 2 .. row 26 all receive the same set)
$values = @(0.6383931775788736, -23.10748969621655, 0.7582052481654122, 0.8998571001689586, 0.8278129647012555, 0.2146654303905574, 14.31124727300474, 0.1887239314917386, 0.06423096092399903, 0.1264774462078688, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656)

# New order of model names for rows 2 through 26
$names = @(
    "model_12_8_0",
    "model_12_8_22",
    "model_12_8_21",
    "model_12_8_20",
    "model_12_8_19",
    "model_12_8_18",
    "model_12_8_17",
    "model_12_8_16",
    "model_12_8_15",
    "model_12_8_14",
    "model_12_8_13",
    "model_12_8_23",
    "model_12_8_12",
    "model_12_8_10",
    "model_12_8_9",
    "model_12_8_8",
    "model_12_8_7",
    "model_12_8_6",
    "model_12_8_5",
    "model_12_8_4",
    "model_12_8_3",
    "model_12_8_2",
    "model_12_8_1",
    "model_12_8_11",
    "model_12_8_24"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $values[$c]
    }
}
